$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper-less approach: this script reproduces the diff by
#  1) inserting a brand new row at row 18 (an "Absent" day, 29.12.2017),
#     which pushes the previous rows 18-20 down to 19-21, and
#  2) appending five brand new rows (22-26) at the bottom of the table for
#     the dates 05.01.2018 .. 11.01.2018.
#
# Dates / time-ranges (e.g. "05.01.2018", "8.40 to 1.00") look like numbers
# to Excel's smart-typing, so every such cell is first switched to the
# "@" (Text) number format, filled in, and then has its formatting restored
# (via copy/PasteSpecial of a same-styled neighbour cell) so the final
# style index matches the rest of the sheet exactly.
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122

# ===========================================================================
# STEP 1 - insert a new row 18 : 13 | Prabha | 29.12.2017 | (blank) | Absent | incomplete
# ===========================================================================
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).RowHeight = 15.75

# Pull cell formatting from row 4 (the existing "Absent/incomplete" row) onto
# the new row's cells - column D is intentionally left untouched (the Insert
# already gave it an empty, styled cell, exactly like the target).
$ws.Range("A4").Copy()
$ws.Range("A18").PasteSpecial($xlPasteFormats)
$ws.Range("B4").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("C4").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("E4").Copy()
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$ws.Range("F4").Copy()
$ws.Range("F18").PasteSpecial($xlPasteFormats)

$ws.Range("A18").Value = 13
$ws.Range("B18").Value = "Prabha"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "29.12.2017"
$ws.Range("E18").Value = "Absent"
$ws.Range("F18").Value = "incomplete"

# Re-apply the formats so the number-format switch above doesn't leave the
# cell on a stray "Text" style - collapses back onto the shared style index.
$ws.Range("A4").Copy()
$ws.Range("A18").PasteSpecial($xlPasteFormats)
$ws.Range("B4").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("C4").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("E4").Copy()
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$ws.Range("F4").Copy()
$ws.Range("F18").PasteSpecial($xlPasteFormats)

# ===========================================================================
# STEP 2 - append rows 22-26 with the new diary entries
# ===========================================================================

function Set-TextCell($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

function Restore-RowFormat($templateRow, $targetRow, $skipD) {
    $cols = @("A","B","C","D","E","F")
    foreach ($col in $cols) {
        if ($skipD -and $col -eq "D") { continue }
        $ws.Range("$col$templateRow").Copy()
        $ws.Range("$col$targetRow").PasteSpecial($xlPasteFormats)
    }
}

# ---- Row 22 : 17 | Prabha | 05.01.2018 | 8.40 to 1.00 | preparing process flow | completed
Restore-RowFormat 21 22 $false
$ws.Range("A22").Value = 17
$ws.Range("B22").Value = "Prabha"
Set-TextCell "C22" "05.01.2018"
Set-TextCell "D22" "8.40 to 1.00"
$ws.Range("E22").Value = "preparing process flow"
$ws.Range("F22").Value = "completed"
$ws.Rows.Item(22).RowHeight = 15.75
Restore-RowFormat 21 22 $false

# ---- Row 23 : 18 | Prabha | 08.01.2018 | 8.30 to 3.45 | connecting bot,JSON with PHP | completed
Restore-RowFormat 21 23 $false
$ws.Range("A23").Value = 18
$ws.Range("B23").Value = "Prabha"
Set-TextCell "C23" "08.01.2018"
Set-TextCell "D23" "8.30 to 3.45"
$ws.Range("E23").Value = "connecting bot,JSON with PHP"
$ws.Range("F23").Value = "completed"
$ws.Rows.Item(23).RowHeight = 15.75
Restore-RowFormat 21 23 $false

# ---- Row 24 : 19 | Prabha | 09.01.2018 | 8.30 to 3.30 | enhancement in flow ,learning php queries  | completed
Restore-RowFormat 21 24 $false
$ws.Range("A24").Value = 19
$ws.Range("B24").Value = "Prabha"
Set-TextCell "C24" "09.01.2018"
Set-TextCell "D24" "8.30 to 3.30"
$ws.Range("E24").Value = "enhancement in flow ,learning php queries "
$ws.Range("F24").Value = "completed"
$ws.Rows.Item(24).RowHeight = 15.75
Restore-RowFormat 21 24 $false

# ---- Row 25 : 20 | Prabha | 10.01.2018 | (blank) | Absent | incomplete
Restore-RowFormat 18 25 $true
$ws.Range("A25").Value = 20
$ws.Range("B25").Value = "Prabha"
Set-TextCell "C25" "10.01.2018"
$ws.Range("E25").Value = "Absent"
$ws.Range("F25").Value = "incomplete"
$ws.Rows.Item(25).RowHeight = 15.75
Restore-RowFormat 18 25 $true

# ---- Row 26 : 21 | Prabha | 11.01.2018 | 8.30 to 4.00 | bot builder API,bot connector-conversational states | completed
Restore-RowFormat 21 26 $false
$ws.Range("A26").Value = 21
$ws.Range("B26").Value = "Prabha"
Set-TextCell "C26" "11.01.2018"
Set-TextCell "D26" "8.30 to 4.00"
$ws.Range("E26").Value = "bot builder API,bot connector-conversational states"
$ws.Range("F26").Value = "completed"
$ws.Rows.Item(26).RowHeight = 15.75
Restore-RowFormat 21 26 $false

# ===========================================================================
# STEP 3 - update the view/selection, mirroring the saved UI state in the diff
# ===========================================================================
$ws.Range("E28").Select()
